$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E are treated as text so that values such as
# "68.517.50", "1.00", "0.0000261" keep their exact original formatting
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "68.517.50"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "3.756.95"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "595.88"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "166.95"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").Value = "3.754.02"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("D11").Value = "6.47"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("D12").Value = "0.448"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "0.0000261"
$ws.Range("E13").Value = "  -6.07%  "
$ws.Range("D14").Value = "36.10"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "4.386.29"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "3.772.30"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "68.487.75"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").Value = "17.93"
$ws.Range("E18").Value = "  -3.96%  "
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "6.99"
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("D21").Value = "10.80"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("D22").Value = "465.18"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "0.696"
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("D24").Value = "84.35"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").Value = "0.0000146"
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("D26").Value = "2.18"
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("D27").Value = "11.92"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  -4.04%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "3.903.82"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").Value = "2.78"
$ws.Range("E31").Value = "  -4.66%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "7.30"
$ws.Range("E32").Value = "  -3.57%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "29.91"
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("D34").Value = "2.16"
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "9.18"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  --%  "
$ws.Range("D37").Value = "3.709.45"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("D38").Value = "0.101"
$ws.Range("E38").Value = "  -3.08%  "
$ws.Range("D39").Value = "3.41"
$ws.Range("E39").Value = "  -7.88%  "
$ws.Range("D40").Value = "0.139"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "5.79"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D45").Value = "0.302"
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("D46").Value = "42.91"
$ws.Range("E46").Value = "  +8.67%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "1.92"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "8.50"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("D49").Value = "46.11"
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("D50").Value = "146.04"
$ws.Range("E50").Value = "  +3.79%  "
$ws.Range("D51").Value = "389.87"
$ws.Range("E51").Value = "  -2.54%  "
